$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 405 (shifts existing rows 405..496 down to 406..497)
$ws.Rows.Item(405).Insert()

# Populate the newly inserted row 405 with the new weekly data record
$ws.Cells.Item(405, 1).Value2 = 9
$ws.Cells.Item(405, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(405, 3).Value2 = "Metropolitana"
$ws.Cells.Item(405, 4).Value2 = 44943
$ws.Cells.Item(405, 5).Value2 = 13
$ws.Cells.Item(405, 6).Value2 = 100112044
$ws.Cells.Item(405, 7).Value2 = "Perejil"
$ws.Cells.Item(405, 8).Value2 = "Sin especificar"
$ws.Cells.Item(405, 9).Value2 = "Primera"
$ws.Cells.Item(405, 10).Value2 = 70
$ws.Cells.Item(405, 11).Value2 = 15000
$ws.Cells.Item(405, 12).Value2 = 17000
$ws.Cells.Item(405, 13).Value2 = 16000
$ws.Cells.Item(405, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(405, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(405, 16).Value2 = 5333
$ws.Cells.Item(405, 17).Value2 = 3
$ws.Cells.Item(405, 18).Value2 = "Hortaliza"
